$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3642143333333334
$ws.Range("H2").Value = 1.092643
$ws.Range("I2").Value = 0.4800482050304226
$ws.Range("J2").Value = 0.4800482050304224
$ws.Range("M2").Value = 166.3936563333333
$ws.Range("N2").Value = 499.180969
$ws.Range("O2").Value = 0.6959913618211631
$ws.Range("P2").Value = 0.7009944564025758
$ws.Range("Q2").Value = 60.60295461234078
$ws.Range("R2").Value = 545.426591511067
$ws.Range("S2").Value = 0.3341094039589287
$ws.Range("T2").Value = 0.3365111305323333

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3642143333333334
$ws.Range("H3").Value = 1.092643
$ws.Range("I3").Value = 0.4800482050304226
$ws.Range("J3").Value = 0.4800482050304224
$ws.Range("O3").Value = 0.2039972194837954
$ws.Range("P3").Value = 0.2054636419703505
$ws.Range("Q3").Value = 17.76291332276167
$ws.Range("R3").Value = 159.866219904855
$ws.Range("S3").Value = 0.09792849904439314
$ws.Range("T3").Value = 0.09863245252688013

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3642143333333334
$ws.Range("H4").Value = 1.092643
$ws.Range("I4").Value = 0.4800482050304226
$ws.Range("J4").Value = 0.4800482050304224
$ws.Range("M4").Value = 7.402863
$ws.Range("N4").Value = 22.208589
$ws.Range("O4").Value = 0.03096469429353687
$ws.Range("P4").Value = 0.03118728224898178
$ws.Range("Q4").Value = 2.696228812303
$ws.Range("R4").Value = 24.266059310727
$ws.Range("S4").Value = 0.01486454591492814
$ws.Range("T4").Value = 0.01497139886340086

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3642143333333334
$ws.Range("H5").Value = 1.092643
$ws.Range("I5").Value = 0.4800482050304226
$ws.Range("J5").Value = 0.4800482050304224
$ws.Range("M5").Value = 5.118919500000001
$ws.Range("N5").Value = 10.237839
$ws.Range("O5").Value = 0.02141141574965316
$ws.Range("P5").Value = 0.0143768870013594
$ws.Range("Q5").Value = 1.8643838530795
$ws.Range("R5").Value = 11.186303118477
$ws.Range("S5").Value = 0.01027851169778112
$ws.Range("T5").Value = 0.006901598798927795

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3642143333333334
$ws.Range("H6").Value = 1.092643
$ws.Range("I6").Value = 0.4800482050304226
$ws.Range("J6").Value = 0.4800482050304224
$ws.Range("M6").Value = 11.38837866666667
$ws.Range("N6").Value = 34.165136
$ws.Range("O6").Value = 0.04763530865185137
$ws.Range("P6").Value = 0.04797773237673265
$ws.Range("Q6").Value = 4.147810743827557
$ws.Range("R6").Value = 37.330296694448
$ws.Range("S6").Value = 0.0228672444143914
$ws.Range("T6").Value = 0.02303162430888049

$ws.Range("G7").Value = 0.3944893333333333
$ws.Range("H7").Value = 1.183468
$ws.Range("I7").Value = 0.5199517949695774
$ws.Range("J7").Value = 0.5199517949695774
$ws.Range("M7").Value = 166.3936563333333
$ws.Range("N7").Value = 499.180969
$ws.Range("O7").Value = 0.6959913618211631
$ws.Range("P7").Value = 0.7009944564025758
$ws.Range("Q7").Value = 65.64052255783244
$ws.Range("R7").Value = 590.764703020492
$ws.Range("S7").Value = 0.3618819578622344
$ws.Range("T7").Value = 0.3644833258702425

$ws.Range("G8").Value = 0.3944893333333333
$ws.Range("H8").Value = 1.183468
$ws.Range("I8").Value = 0.5199517949695774
$ws.Range("J8").Value = 0.5199517949695774
$ws.Range("O8").Value = 0.2039972194837954
$ws.Range("P8").Value = 0.2054636419703505
$ws.Range("Q8").Value = 19.23944005888667
$ws.Range("R8").Value = 173.15496052998
$ws.Range("S8").Value = 0.1060687204394023
$ws.Range("T8").Value = 0.1068311894434704

$ws.Range("G9").Value = 0.3944893333333333
$ws.Range("H9").Value = 1.183468
$ws.Range("I9").Value = 0.5199517949695774
$ws.Range("J9").Value = 0.5199517949695774
$ws.Range("M9").Value = 7.402863
$ws.Range("N9").Value = 22.208589
$ws.Range("O9").Value = 0.03096469429353687
$ws.Range("P9").Value = 0.03118728224898178
$ws.Range("Q9").Value = 2.920350489628
$ws.Range("R9").Value = 26.283154406652
$ws.Range("S9").Value = 0.01610014837860873
$ws.Range("T9").Value = 0.01621588338558092

$ws.Range("G10").Value = 0.3944893333333333
$ws.Range("H10").Value = 1.183468
$ws.Range("I10").Value = 0.5199517949695774
$ws.Range("J10").Value = 0.5199517949695774
$ws.Range("M10").Value = 5.118919500000001
$ws.Range("N10").Value = 10.237839
$ws.Range("O10").Value = 0.02141141574965316
$ws.Range("P10").Value = 0.0143768870013594
$ws.Range("Q10").Value = 2.019359140942
$ws.Range("R10").Value = 12.116154845652
$ws.Range("S10").Value = 0.01113290405187204
$ws.Range("T10").Value = 0.007475288202431608

$ws.Range("G11").Value = 0.3944893333333333
$ws.Range("H11").Value = 1.183468
$ws.Range("I11").Value = 0.5199517949695774
$ws.Range("J11").Value = 0.5199517949695774
$ws.Range("M11").Value = 11.38837866666667
$ws.Range("N11").Value = 34.165136
$ws.Range("O11").Value = 0.04763530865185137
$ws.Range("P11").Value = 0.04797773237673265
$ws.Range("Q11").Value = 4.492593907960889
$ws.Range("R11").Value = 40.433345171648
$ws.Range("S11").Value = 0.02476806423745996
$ws.Range("T11").Value = 0.02494610806785215

